$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: updated timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now set
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 (was "Contact" / "No display for ContactDetail") becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely
$ws.Rows.Item(11).Delete()

# --- Sheet "Elements" ---
$ws2 = $wb.Worksheets.Item("Elements")

# Row 2 (root Extension element): Short & Definition updated to be specific
$ws2.Range("K2").Value = "Funding Arrangement"
$ws2.Range("L2").Value = "Code for the funding arrangement"
